$d = $word.ActiveDocument

$d.Content.Find.Execute("Based on the two indices you found for", $true, $false, $false, $false, $false, $true, 1, $false, "Based on the VIS indices you found out for", 2)
$d.Content.Find.Execute("homework 1 and the mosaics", $true, $false, $false, $false, $false, $true, 1, $false, "homework 2 and the mosaics", 2)
$d.Content.Find.Execute("the correlation between them and those used in class", $true, $false, $false, $false, $false, $true, 1, $false, "the correlation between yours and the ones used in class", 2)
